# "4.29 Advance the story line"
#
# - Row 10 ("End of story") is retired; a new row 7 ("Goto" / "StoryScript2")
#   is added right after the existing dialogue rows so the script can jump
#   back into another story chunk.
# - Four new helper columns (M:P = LastBGImage / LastBGM / LastX1 / LastX2)
#   are added. For every dialogue row they carry forward the most recent
#   non-blank BGImage / BGM / x1 / x2 value from the row above, so a "Goto"
#   can restore whatever background/position was last in effect.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old "End of story" row (A10) with the new Goto/StoryScript2 row
# (A7/B7). Doing this before the new header cells keeps the shared-string
# insertion order (Goto right after "appearAt", then the Last* headers,
# then StoryScript2) the same as a real edit session would produce.
$ws.Cells.Item(10, 1).ClearContents()
$ws.Range("A7").Value = "Goto"

# New header cells for the carried-forward helper columns.
$ws.Range("M1").Value = "LastBGImage"
$ws.Range("N1").Value = "LastBGM"
$ws.Range("O1").Value = "LastX1"
$ws.Range("P1").Value = "LastX2"

$ws.Range("B7").Value = "StoryScript2"

# First row of the helper columns falls back to row 2 (the very first
# dialogue row) when a value is blank.
$ws.Range("M3").Formula = '=IF(E2<>"",E2,M2)'
$ws.Range("N3").Formula = '=IF(F2<>"",F2,N2)'
$ws.Range("O3").Formula = '=IF(H2<>"",H2,O2)'
$ws.Range("P3").Formula = '=IF(K2<>"",K2,P2)'

# Rows 4:6 repeat the same "carry forward the previous row" pattern, filled
# down as one shared formula per column.
$ws.Range("M4:M6").Formula = '=IF(E3<>"",E3,M3)'
$ws.Range("N4:N6").Formula = '=IF(F3<>"",F3,N3)'
$ws.Range("O4:O6").Formula = '=IF(H3<>"",H3,O3)'
$ws.Range("P4:P6").Formula = '=IF(K3<>"",K3,P3)'

# Match the row's rendered (single-line) height with the other plain rows.
$ws.Rows.Item(7).RowHeight = 17

# Leave the selection where the author ended up after typing the new row.
$ws.Range("B8").Select()
